$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

$ws.Range("H41").Value = 1320.875
$ws.Range("I41").Value = 2337
$ws.Range("J41").Value = 304.75
$ws.Range("K41").Value = 2337
$ws.Range("L41").Value = 304.75
$ws.Range("M41").Value = -1897
$ws.Range("N41").Value = -1184.75

$ws.Range("H51").Value = 7832.857
$ws.Range("I51").Value = 27745
$ws.Range("J51").Value = 5736.8423
$ws.Range("K51").Value = 27745
$ws.Range("L51").Value = 5736.8423
$ws.Range("M51").Value = -27261

$ws.Range("H53").Value = 406.63635
$ws.Range("I53").Value = 457.5
$ws.Range("J53").Value = 387.5625
$ws.Range("K53").Value = 457.5
$ws.Range("L53").Value = 387.5625
$ws.Range("M53").Value = 179.5
$ws.Range("N53").Value = -1661.5625

$ws.Range("H55").Value = 3803.8
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3803.8
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 3803.8
$ws.Range("N55").Value = -4231.8

$ws.Range("H92").Value = 502.8
$ws.Range("I92").Value = 82.59999999999999
$ws.Range("J92").Value = 1343.2
$ws.Range("K92").Value = 82.59999999999999
$ws.Range("L92").Value = 1343.2
$ws.Range("M92").Value = 1165.4
$ws.Range("N92").Value = -3839.2

$ws.Range("H100").Value = 9401.071
$ws.Range("I100").Value = 5327.1665
$ws.Range("J100").Value = 12456.5
$ws.Range("K100").Value = 5327.1665
$ws.Range("L100").Value = 12456.5
$ws.Range("M100").Value = -4786.1665
$ws.Range("N100").Value = -13538.5

$ws.Range("H107").Value = 524.4
$ws.Range("I107").Value = 498.57144
$ws.Range("J107").Value = 886
$ws.Range("K107").Value = 498.57144
$ws.Range("L107").Value = 886
$ws.Range("M107").Value = 1421.42856
$ws.Range("N107").Value = -4726

$ws.Range("H116").Value = 9722.654
$ws.Range("I116").Value = 5497.5713
$ws.Range("J116").Value = 11279.263
$ws.Range("K116").Value = 5497.5713
$ws.Range("L116").Value = 11279.263
$ws.Range("M116").Value = -2055.5713
$ws.Range("N116").Value = -18163.263

$ws.Range("H127").Value = 2696.1667
$ws.Range("I127").Value = 1549.25
$ws.Range("J127").Value = 4990
$ws.Range("K127").Value = 4647.75
$ws.Range("L127").Value = 14970
$ws.Range("M127").Value = 312.25

$ws.Range("H132").Value = 4120.489
$ws.Range("I132").Value = 2536.861
$ws.Range("J132").Value = 10455
$ws.Range("K132").Value = 7610.583
$ws.Range("L132").Value = 31365
$ws.Range("M132").Value = -5080.583
$ws.Range("N132").Value = -36425

$ws.Range("H137").Value = 3527
$ws.Range("I137").Value = 2276.9092
$ws.Range("J137").Value = 6277.2
$ws.Range("K137").Value = 6830.7276
$ws.Range("L137").Value = 18831.6
$ws.Range("M137").Value = -4280.7276
$ws.Range("N137").Value = -23931.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1625250
$ws.Range("I8").Value = 1857142.9
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 1857142.9
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = -1856998.9
$ws.Range("N8").Value = -2288

$ws.Range("H10").Value = 6969
$ws.Range("I10").Value = 6969
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 6969
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -6799
$ws.Range("N10").Value = ""

$ws.Range("H12").Value = 1574.875
$ws.Range("I12").Value = 999.5
$ws.Range("J12").Value = 1766.6666
$ws.Range("K12").Value = 999.5
$ws.Range("L12").Value = 1766.6666
$ws.Range("M12").Value = -826.5
$ws.Range("N12").Value = -2112.6666

$ws.Range("H13").Value = 1166666.6
$ws.Range("I13").Value = 1166666.6
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1166666.6
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1166522.6
$ws.Range("N13").Value = ""

$ws.Range("H32").Value = 4565.735
$ws.Range("I32").Value = 3982.689
$ws.Range("J32").Value = 11125
$ws.Range("K32").Value = 3982.689
$ws.Range("L32").Value = 11125
$ws.Range("M32").Value = -3695.689

$ws.Range("H43").Value = 32170.5
$ws.Range("I43").Value = 39195
$ws.Range("J43").Value = 30414.375
$ws.Range("K43").Value = 39195
$ws.Range("L43").Value = 30414.375
$ws.Range("M43").Value = -38882
$ws.Range("N43").Value = -31040.375

$ws.Range("H61").Value = 12224648
$ws.Range("I61").Value = 15002691
$ws.Range("J61").Value = 2501496.8
$ws.Range("K61").Value = 15002691
$ws.Range("L61").Value = 2501496.8
$ws.Range("M61").Value = -15002479

$ws.Range("H74").Value = 3566.1
$ws.Range("I74").Value = 3237
$ws.Range("J74").Value = 5431
$ws.Range("K74").Value = 3237
$ws.Range("L74").Value = 5431
$ws.Range("M74").Value = -2363

$ws.Range("H77").Value = 3566.1
$ws.Range("I77").Value = 3237
$ws.Range("J77").Value = 5431
$ws.Range("K77").Value = 16185
$ws.Range("L77").Value = 27155
$ws.Range("M77").Value = -11817

$ws.Range("H112").Value = 11380
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 11380
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 11380
$ws.Range("N112").Value = -14334

$ws.Range("H132").Value = 2276210.5
$ws.Range("I132").Value = 3642.853
$ws.Range("J132").Value = 10002941
$ws.Range("K132").Value = 10928.559
$ws.Range("L132").Value = 30008823
$ws.Range("M132").Value = -8398.559000000001
$ws.Range("N132").Value = -30013883

$ws.Range("H136").Value = 12224648
$ws.Range("I136").Value = 15002691
$ws.Range("J136").Value = 2501496.8
$ws.Range("K136").Value = 45008073
$ws.Range("L136").Value = 7504490.399999999
$ws.Range("M136").Value = -45005523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 600
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -940

$ws.Range("H20").Value = 8780.77
$ws.Range("I20").Value = 10986.2
$ws.Range("J20").Value = 1429.3334
$ws.Range("K20").Value = 10986.2
$ws.Range("L20").Value = 1429.3334
$ws.Range("M20").Value = -10739.2
$ws.Range("N20").Value = -1923.3334

$ws.Range("H110").Value = 99990
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 99990
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 99990
$ws.Range("N110").Value = -108170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52634772
$ws.Range("I31").Value = 83336420
$ws.Range("J31").Value = 3386.1428
$ws.Range("K31").Value = 83336420
$ws.Range("L31").Value = 3386.1428
$ws.Range("M31").Value = -83336125
$ws.Range("N31").Value = -3976.1428

$ws.Range("H34").Value = 52634772
$ws.Range("I34").Value = 83336420
$ws.Range("J34").Value = 3386.1428
$ws.Range("K34").Value = 83336420
$ws.Range("L34").Value = 3386.1428
$ws.Range("M34").Value = -83336218
$ws.Range("N34").Value = -3790.1428

$ws.Range("H59").Value = 50666
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 60999
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 60999
$ws.Range("M59").Value = -28855
$ws.Range("N59").Value = -63289

$ws.Range("H122").Value = 3695.85
$ws.Range("I122").Value = 3344.4614
$ws.Range("J122").Value = 4348.4287
$ws.Range("K122").Value = 10033.3842
$ws.Range("L122").Value = 13045.2861
$ws.Range("M122").Value = -7583.3842

$ws.Range("H134").Value = 1815.9166
$ws.Range("I134").Value = 1781.091
$ws.Range("J134").Value = 2199
$ws.Range("K134").Value = 5343.272999999999
$ws.Range("L134").Value = 6597
$ws.Range("M134").Value = -2808.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 18395.75
$ws.Range("I94").Value = 250
$ws.Range("J94").Value = 24444.334
$ws.Range("K94").Value = 750
$ws.Range("L94").Value = 73333.00199999999
$ws.Range("M94").Value = -74
$ws.Range("N94").Value = -74685.00199999999

$ws.Range("H107").Value = 4337225.5
$ws.Range("I107").Value = 3060.2
$ws.Range("J107").Value = 5691652.5
$ws.Range("K107").Value = 9180.599999999999
$ws.Range("L107").Value = 17074957.5
$ws.Range("M107").Value = -7260.599999999999
$ws.Range("N107").Value = -17078797.5

$ws.Range("H129").Value = 16671932
$ws.Range("I129").Value = 22729784
$ws.Range("J129").Value = 12839
$ws.Range("K129").Value = 68189352
$ws.Range("L129").Value = 38517
$ws.Range("M129").Value = -68184352
$ws.Range("N129").Value = -48517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2808.2856
$ws.Range("I102").Value = 2698.75
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 2698.75
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = -1076.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4553.643
$ws.Range("I7").Value = 3980.923
$ws.Range("J7").Value = 11999
$ws.Range("K7").Value = 3980.923
$ws.Range("L7").Value = 11999
$ws.Range("M7").Value = -3868.923
$ws.Range("N7").Value = -12223

$ws.Range("H82").Value = 4703.0557
$ws.Range("I82").Value = 2729.5557
$ws.Range("J82").Value = 6676.5557
$ws.Range("K82").Value = 2729.5557
$ws.Range("L82").Value = 6676.5557
$ws.Range("M82").Value = -2368.5557
$ws.Range("N82").Value = -7398.5557

$ws.Range("H85").Value = 4703.0557
$ws.Range("I85").Value = 2729.5557
$ws.Range("J85").Value = 6676.5557
$ws.Range("K85").Value = 2729.5557
$ws.Range("L85").Value = 6676.5557
$ws.Range("M85").Value = -1481.5557
$ws.Range("N85").Value = -9172.555700000001

$ws.Range("H110").Value = 71666
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 71666
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 71666
$ws.Range("N110").Value = -79846

$ws.Range("H126").Value = 4553.643
$ws.Range("I126").Value = 3980.923
$ws.Range("J126").Value = 11999
$ws.Range("K126").Value = 11942.769
$ws.Range("L126").Value = 35997
$ws.Range("M126").Value = -9472.769
$ws.Range("N126").Value = -40937

$ws.Range("H136").Value = 2428.9395
$ws.Range("I136").Value = 2460.5417
$ws.Range("J136").Value = 2344.6667
$ws.Range("K136").Value = 7381.625100000001
$ws.Range("L136").Value = 7034.000100000001
$ws.Range("M136").Value = -4831.625100000001
$ws.Range("N136").Value = -12134.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16421.857
$ws.Range("I41").Value = 26995
$ws.Range("J41").Value = 14659.667
$ws.Range("K41").Value = 26995
$ws.Range("L41").Value = 14659.667
$ws.Range("M41").Value = -26605
$ws.Range("N41").Value = -15439.667

$ws.Range("H136").Value = 278197.3
$ws.Range("I136").Value = 9173.379000000001
$ws.Range("J136").Value = 1253409.1
$ws.Range("K136").Value = 27520.137
$ws.Range("L136").Value = 3760227.3
$ws.Range("M136").Value = -24970.137
